# "modify map to make prd"
# Replace the seven residential-project names in column A with the
# mining-enterprise names used for the PRD map, shift the longitude
# (column C) values from the 108/109 range into the 107 range, widen
# column A to fit the new (longer) labels, and update the saved
# selection / scroll position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- column A: project/enterprise names -----------------------------------
$ws.Range("A2").Value = "铅锌冶炼企业一"
$ws.Range("A3").Value = "铅锌冶炼企业二"
$ws.Range("A4").Value = "铜冶炼企业一"
$ws.Range("A5").Value = "金矿采选企业一"
$ws.Range("A6").Value = "金矿采选企业二"
$ws.Range("A7").Value = "钨钼冶炼企业一"
$ws.Range("A8").Value = "钨钼冶炼企业二"

# --- column C: longitude values shift from 108/109.x to 107.x -------------
$ws.Range("C2").Value = 107.86695899999999
$ws.Range("C3").Value = 107.870649
$ws.Range("C4").Value = 107.15137300000001
$ws.Range("C5").Value = 107.19124100000001
$ws.Range("C6").Value = 107.19767899999999
$ws.Range("C7").Value = 107.04342699999999
$ws.Range("C8").Value = 107.876572

# --- column A width: widen to fit the new labels ---------------------------
$ws.Columns.Item(1).ColumnWidth = 14.25

# --- view state: scroll/selection moved to A13, no pinned top-left cell ----
$ws.Activate() | Out-Null
$ws.Range("A13").Select() | Out-Null

# --- workbook window position reset to 0,0 ---------------------------------
try { $excel.Left = 0 } catch {}
try { $excel.Top = 0 } catch {}
try { $wb.Windows.Item(1).Left = 0 } catch {}
try { $wb.Windows.Item(1).Top = 0 } catch {}
